$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.03407883644104
$ws.Range("B1").Value = 2.132187366485596
$ws.Range("C1").Value = 2.534010648727417
$ws.Range("D1").Value = 2.761887550354004
$ws.Range("E1").Value = 1.415645599365234
